$d = $word.ActiveDocument

# --- Change 1: merge the "An Employer ... username " run with the "ad" run
# (dropping the spell-check proofErr markers around "add"), while leaving the
# trailing "d" run and the " password." run untouched. -----------------------
$p9 = $d.Paragraphs(9)
$r9 = $p9.Range
# Exclude the trailing paragraph mark from the range we rewrite.
$bodyRange = $d.Range($r9.Start, $r9.End - 1)
$newRunsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
  + '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr>' `
  + '<w:t>An Employer would receive an email from the TDA Admin containing their username ad</w:t></w:r>' `
  + '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>d</w:t></w:r>' `
  + '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> password.</w:t></w:r>' `
  + '</w:p>'
$bodyRange.InsertXML($newRunsXml)

# --- Change 2: add a new bullet paragraph after "email link." -------------
$p13 = $d.Paragraphs(13)
$p13.Range.InsertParagraphAfter()
$p14 = $d.Paragraphs(14)
$p14.Range.Text = "Ideally would have profile picture out but due to time constraints this would future workings."

Write-Output "done"
